# Fix up a few mis-placed "Description" / "Value / Tag" entries on the
# "MDI Tags and fields" worksheet (Table1, columns A=Name/Label,
# B=Description, C=Value/Tag).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 44 "Bits in" - value was already present (C44='bitsin') but the
# description was blank. Add the missing description.
$ws.Range("B44").Value = "In the context of the Impacted Host."

# Row 45 "Bits out" - the value ('bitsout') had been typed into the
# Description column (B) instead of the Value / Tag column (C). Move it
# over and supply the correct description text.
$ws.Range("C45").Value = "bitsout"
$ws.Range("B45").Value = "In the context of the Impacted Host."

# Row 86 "Vendor Message ID" - the value ('vmid') had been appended onto
# the end of the Description text instead of its own column. Split them
# back apart.
$ws.Range("B86").Value = "Specific vendor for the log used to describe a type of event."
$ws.Range("C86").Value = "vmid"
